# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$geral = $wb.Worksheets.Item("Geral")

# Update Rodada 1 scores on the 'Geral' sheet (column B)
$geral.Range("B2").Value = 68.06
$geral.Range("B3").Value = 54.1
$geral.Range("B4").Value = 38.26
$geral.Range("B5").Value = 74.06
$geral.Range("B6").Value = 62.26
$geral.Range("B8").Value = 61.56
$geral.Range("B9").Value = 73.76000000000001
$geral.Range("B10").Value = 72.86
$geral.Range("B11").Value = 60.2
$geral.Range("B12").Value = 72.7
$geral.Range("B13").Value = 53.66
$geral.Range("B14").Value = 56.96
$geral.Range("B15").Value = 37.82
$geral.Range("B16").Value = 56.26
$geral.Range("B17").Value = 63.9
$geral.Range("B19").Value = 54.16
$geral.Range("B20").Value = 40.6
$geral.Range("B21").Value = 57.45
$geral.Range("B22").Value = 57.6
$geral.Range("B23").Value = 54.6
$geral.Range("B24").Value = 47.86
$geral.Range("B25").Value = 47.86
$geral.Range("B26").Value = 62.56
$geral.Range("B27").Value = 57.66
$geral.Range("B29").Value = 73.95999999999999
$geral.Range("B30").Value = 51.56
$geral.Range("B31").Value = 49.76
$geral.Range("B32").Value = 43.46
$geral.Range("B33").Value = 84.26000000000001
$geral.Range("B34").Value = 64.56
$geral.Range("B35").Value = 43.51
$geral.Range("B36").Value = 57.46
$geral.Range("B37").Value = 64.56
$geral.Range("B38").Value = 59.25
$geral.Range("B39").Value = 84.86
$geral.Range("B40").Value = 60.16
$geral.Range("B41").Value = 49
$geral.Range("B42").Value = 61.66
$geral.Range("B43").Value = 42.96
$geral.Range("B44").Value = 54.1
$geral.Range("B45").Value = 61.96
$geral.Range("B46").Value = 64.2
$geral.Range("B48").Value = 58.96
$geral.Range("B49").Value = 70
$geral.Range("B50").Value = 68.06
$geral.Range("B51").Value = 51.26
$geral.Range("B52").Value = 72.16
$geral.Range("B53").Value = 55.9

$janeiro = $wb.Worksheets.Item("Mes - Janeiro")

# Rewrite the January ranking table (columns A and B), re-sorted by score desc
$janeiro.Cells.Item(2, 1).Value = "S.E.R. GRILLO"
$janeiro.Cells.Item(2, 2).Value = 84.86
$janeiro.Cells.Item(3, 1).Value = "Paulo Virgili FC"
$janeiro.Cells.Item(3, 2).Value = 84.26000000000001
$janeiro.Cells.Item(4, 1).Value = "Bandoleros FCS"
$janeiro.Cells.Item(4, 2).Value = 74.06
$janeiro.Cells.Item(5, 1).Value = "Mau Humor F.C."
$janeiro.Cells.Item(5, 2).Value = 73.95999999999999
$janeiro.Cells.Item(6, 1).Value = "CARTOLEIRO DO VALLE PRO26.5"
$janeiro.Cells.Item(6, 2).Value = 73.76000000000001
$janeiro.Cells.Item(7, 1).Value = "dasdoresfc"
$janeiro.Cells.Item(7, 2).Value = 72.86
$janeiro.Cells.Item(8, 1).Value = "Dom Camillo68"
$janeiro.Cells.Item(8, 2).Value = 72.7
$janeiro.Cells.Item(9, 1).Value = "Time do S.A.P.O"
$janeiro.Cells.Item(9, 2).Value = 72.16
$janeiro.Cells.Item(10, 1).Value = "C.A. Charru@"
$janeiro.Cells.Item(10, 2).Value = 71.70999999999999
$janeiro.Cells.Item(11, 1).Value = "teves_futsal20 f.c"
$janeiro.Cells.Item(11, 2).Value = 70
$janeiro.Cells.Item(12, 1).Value = "Texas Club 2026"
$janeiro.Cells.Item(12, 2).Value = 68.06
$janeiro.Cells.Item(13, 1).Value = "A Lenda Super Vasco F.c"
$janeiro.Cells.Item(13, 2).Value = 68.06
$janeiro.Cells.Item(14, 1).Value = "Tatols Beants F.C"
$janeiro.Cells.Item(14, 2).Value = 66.86
$janeiro.Cells.Item(15, 1).Value = "PUXE FC"
$janeiro.Cells.Item(15, 2).Value = 64.56
$janeiro.Cells.Item(16, 1).Value = "Pity10"
$janeiro.Cells.Item(16, 2).Value = 64.56
$janeiro.Cells.Item(17, 1).Value = "TATITTA FC"
$janeiro.Cells.Item(17, 2).Value = 64.2
$janeiro.Cells.Item(18, 1).Value = "Fedato Futebol Clube"
$janeiro.Cells.Item(18, 2).Value = 63.9
$janeiro.Cells.Item(19, 1).Value = "lsauer fc"
$janeiro.Cells.Item(19, 2).Value = 62.56
$janeiro.Cells.Item(20, 1).Value = "BordonFC04"
$janeiro.Cells.Item(20, 2).Value = 62.26
$janeiro.Cells.Item(21, 1).Value = "Tabajara de Inhaua PB1"
$janeiro.Cells.Item(21, 2).Value = 61.96
$janeiro.Cells.Item(22, 1).Value = "Sport Clube PAIM"
$janeiro.Cells.Item(22, 2).Value = 61.66
$janeiro.Cells.Item(23, 1).Value = "cartola scheuer17"
$janeiro.Cells.Item(23, 2).Value = 61.56
$janeiro.Cells.Item(24, 1).Value = "DM Studio"
$janeiro.Cells.Item(24, 2).Value = 60.2
$janeiro.Cells.Item(25, 1).Value = "seralex"
$janeiro.Cells.Item(25, 2).Value = 60.16
$janeiro.Cells.Item(26, 1).Value = "Rolo Compressor ZN"
$janeiro.Cells.Item(26, 2).Value = 59.25
$janeiro.Cells.Item(27, 1).Value = "TEAM LOPES 99"
$janeiro.Cells.Item(27, 2).Value = 58.96
$janeiro.Cells.Item(28, 1).Value = "MAFRA MARTINS FC"
$janeiro.Cells.Item(28, 2).Value = 58.51
$janeiro.Cells.Item(29, 1).Value = "Luis lemes inter"
$janeiro.Cells.Item(29, 2).Value = 57.66
$janeiro.Cells.Item(30, 1).Value = "Grêmio imortal 37"
$janeiro.Cells.Item(30, 2).Value = 57.6
$janeiro.Cells.Item(31, 1).Value = "Profit Soccer"
$janeiro.Cells.Item(31, 2).Value = 57.46
$janeiro.Cells.Item(32, 1).Value = "Gremiomaniasm"
$janeiro.Cells.Item(32, 2).Value = 57.45
$janeiro.Cells.Item(33, 1).Value = "FBC Colorado II"
$janeiro.Cells.Item(33, 2).Value = 56.96
$janeiro.Cells.Item(34, 1).Value = "FC Los Castilho"
$janeiro.Cells.Item(34, 2).Value = 56.26
$janeiro.Cells.Item(35, 1).Value = "VASCO MARTINS FC"
$janeiro.Cells.Item(35, 2).Value = 55.9
$janeiro.Cells.Item(36, 1).Value = "JUV. KP"
$janeiro.Cells.Item(36, 2).Value = 54.6
$janeiro.Cells.Item(37, 1).Value = "FÚRIA LEON"
$janeiro.Cells.Item(37, 2).Value = 54.16
$janeiro.Cells.Item(38, 1).Value = "SUPER VASCÃO F.C"
$janeiro.Cells.Item(38, 2).Value = 54.1
$janeiro.Cells.Item(39, 1).Value = "A Lenda Super Vascão f.c"
$janeiro.Cells.Item(39, 2).Value = 54.1
$janeiro.Cells.Item(40, 1).Value = "FBC Colorado"
$janeiro.Cells.Item(40, 2).Value = 53.66
$janeiro.Cells.Item(41, 1).Value = "mercearia Estrela"
$janeiro.Cells.Item(41, 2).Value = 51.56
$janeiro.Cells.Item(42, 1).Value = "TIGRE LEON"
$janeiro.Cells.Item(42, 2).Value = 51.26
$janeiro.Cells.Item(43, 1).Value = "Máquina Laranjja"
$janeiro.Cells.Item(43, 2).Value = 49.76
$janeiro.Cells.Item(44, 1).Value = "SERGRILLO"
$janeiro.Cells.Item(44, 2).Value = 49
$janeiro.Cells.Item(45, 1).Value = "JV5 Tricolor Gaúcho"
$janeiro.Cells.Item(45, 2).Value = 47.86
$janeiro.Cells.Item(46, 1).Value = "LISI GREMISTA"
$janeiro.Cells.Item(46, 2).Value = 47.86
$janeiro.Cells.Item(47, 1).Value = "FIGUEIRA DA ILHA"
$janeiro.Cells.Item(47, 2).Value = 44.06
$janeiro.Cells.Item(48, 1).Value = "pra sempre imortal fc"
$janeiro.Cells.Item(48, 2).Value = 43.51
$janeiro.Cells.Item(49, 1).Value = "NaoVaiDescer!"
$janeiro.Cells.Item(49, 2).Value = 43.46
$janeiro.Cells.Item(50, 1).Value = "Super Vasco f.c"
$janeiro.Cells.Item(50, 2).Value = 42.96
$janeiro.Cells.Item(51, 1).Value = "Gig@ntte"
$janeiro.Cells.Item(51, 2).Value = 40.6
$janeiro.Cells.Item(52, 1).Value = "AZURRA82"
$janeiro.Cells.Item(52, 2).Value = 38.26
$janeiro.Cells.Item(53, 1).Value = "FC castelo Branco 2"
$janeiro.Cells.Item(53, 2).Value = 37.82
